# MENT-169: Create New Question Categories
#
# Inserts a new "Área Clinica" / "Oficial de Pediatria" row just above the
# existing "Monitoria e Avaliação" block (which currently occupies rows
# 98-100), pushing that block down by one row (98-100 -> 99-101) and
# renumbering the "Nr" column accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the three existing "Monitoria e Avaliação" rows down and make room
# for the new category row at row 98.
$ws.Rows.Item(98).Insert()

# New row: Área Clinica / Oficial de Pediatria
$ws.Cells.Item(98, 1).Value = 96
$ws.Cells.Item(98, 2).Value = "CLINICAL_AREA"
$ws.Cells.Item(98, 3).Value = "Área Clinica"
$ws.Cells.Item(98, 4).Value = "Oficial de Pediatria"

# Renumber the "Nr" column for the rows that shifted down.
$ws.Cells.Item(99, 1).Value = 97
$ws.Cells.Item(100, 1).Value = 98
$ws.Cells.Item(101, 1).Value = 99

# Restore the selection to match the post-edit view.
$ws.Range("B20").Select() | Out-Null
